# Weekly update: insert a new price record as row 37 for
# "Vega Monumental Concepción - Haba", shifting existing rows 37-81 down to 38-82.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 37 (pushes old rows 37..81 down to 38..82,
# carrying over formatting from the row above, same as Excel's UI "Insert").
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Range("A37").Value = 11
$ws.Range("B37").Value = "Vega Monumental Concepción"
$ws.Range("C37").Value = "Bíobío"
$ws.Range("D37").Value = 45272
$ws.Range("E37").Value = 8
$ws.Range("F37").Value = 100112026
$ws.Range("G37").Value = "Haba"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 120
$ws.Range("K37").Value = 13000
$ws.Range("L37").Value = 13000
$ws.Range("M37").Value = 13000
$ws.Range("N37").Value = "`$/saco 25 kilos"
$ws.Range("O37").Value = "Región del Maule"
$ws.Range("P37").Value = 520
$ws.Range("Q37").Value = 25
$ws.Range("R37").Value = "Hortaliza"
